$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill cells in the same order the strings were originally entered so the
# shared-strings table comes out in the matching sequence.
$ws.Range("A1").Value = "`${test_case_name}"
$ws.Range("A2").Value = "add_valid"
$ws.Range("B1").Value = "`${year}"
$ws.Range("B2").Value = "2021-2025"
$ws.Range("C1").Value = "`${bussiness}"
$ws.Range("D1").Value = "`${service1}"
$ws.Range("C2").Value = "training"
$ws.Range("D2").Value = "testing automation"
$ws.Range("E1").Value = "`${business2}"
$ws.Range("F1").Value = "`${service2}"

$ws.Range("A3").Value = "delete_valid"
$ws.Range("B3").Value = "2021-2026"
$ws.Range("C3").Value = "csr"
$ws.Range("D3").Value = "test servies for deletion"

$ws.Range("A4").Value = "empty_year"
$ws.Range("C4").Value = "csr"
$ws.Range("D4").Value = "sample test service"

$ws.Range("A5").Value = "add_new_service"
$ws.Range("B5").Value = "2021-2026"
$ws.Range("C5").Value = "b2b"
$ws.Range("D5").Value = "testing automation"
$ws.Range("F5").Value = "new service"

$ws.Range("A6").Value = "del_new_service"
$ws.Range("B6").Value = 2024
$ws.Range("C6").Value = "crc"
$ws.Range("D6").Value = "fullstack"
$ws.Range("F6").Value = "deleted service"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 20.109375
$ws.Columns.Item(3).ColumnWidth = 13.109375
$ws.Columns.Item(4).ColumnWidth = 17.21875
$ws.Columns.Item(5).ColumnWidth = 14.6640625
$ws.Columns.Item(6).ColumnWidth = 13.6640625

# Selection
$ws.Range("A6").Select()
